$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title 1: "A slide" -> split "A " into "A" + " " runs ---
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Characters(1, 1).Text = "A"

# --- TextBox 3: "Just an image on this side" ---
# Split each "word " run into separate "word" and " " runs
# (leaving the final "side" run untouched).
$box = $s.Shapes.Item(4)
$boxRange = $box.TextFrame.TextRange
$boxRange.Characters(1, 4).Text = "Just"
$boxRange.Characters(6, 2).Text = "an"
$boxRange.Characters(9, 5).Text = "image"
$boxRange.Characters(15, 2).Text = "on"
$boxRange.Characters(18, 4).Text = "this"
